$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top for "Polygon Window - Surfing on Sine Waves".
$ws.Rows.Item(1).Insert() | Out-Null
$ws.Range("B1").Value = "Polygon Window"
$ws.Range("C1").Value = "Surfing on Sine Waves"
$ws.Range("D1").Value = "Warp Records"
$ws.Range("E1").Value = 1993
$ws.Range("F1").Value = "ELETRONICA"
$ws.Range("G1").Value = "CD"

# Insert another new row above it for "Various - Soulfood".
$ws.Rows.Item(1).Insert() | Out-Null
$ws.Range("B1").Value = "Various"
$ws.Range("C1").Value = "Soulfood"
$ws.Range("D1").Value = "Cookin Records"
$ws.Range("E1").Value = 1999
$ws.Range("F1").Value = "ELETRONICA"
$ws.Range("G1").Value = "CD"

# Assign catalogue numbers now that the final row order is set.
$ws.Range("A1").Value = "c86"
$ws.Range("A2").Value = "c87"

# Rename the genre labels across the whole genre column (whole-cell match,
# so "ELETRONICA" only touches pure-genre cells, not "JAZZ, ELETRONICA").
$genreCol = $ws.Columns.Item(6)
$genreCol.Replace("ELETRONICA", "ELETRONIC", 1) | Out-Null
$genreCol.Replace("JAZZ, ELETRONICA", "JAZZ, ELETRONIC", 1) | Out-Null

# Reset the view so the top-left cell is back at A1 and F1 is selected.
$ws.Range("F1").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
